$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48
$ws.Range("W2").Value = 7
$ws.Range("AC2").Value = 6.5
$ws.Range("AF2").Value = 67
$ws.Range("AY2").Value = 29
$ws.Range("J4").Value = 2.75
$ws.Range("K4").Value = 2.25
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 4
$ws.Range("Q4").Value = 1.75
$ws.Range("R4").Value = 2.05
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 3.25
$ws.Range("U4").Value = 1.62
$ws.Range("V4").Value = 2.2
$ws.Range("W4").Value = 9.5
$ws.Range("X4").Value = 12
$ws.Range("AC4").Value = 12
$ws.Range("AG4").Value = 12
$ws.Range("AT4").Value = 3.25
$ws.Range("G5").Value = 1.62
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 6.25
$ws.Range("J5").Value = 2.3
$ws.Range("L5").Value = 6.5
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 6.5
$ws.Range("Z5").Value = 11
$ws.Range("AD5").Value = 7
$ws.Range("AG5").Value = 12
$ws.Range("AH5").Value = 29
$ws.Range("AI5").Value = 21
$ws.Range("AL5").Value = 67
$ws.Range("AN5").Value = 3.4
$ws.Range("AO5").Value = 8.5
$ws.Range("AP5").Value = 23
$ws.Range("AQ5").Value = 29
$ws.Range("AR5").Value = 51
$ws.Range("AU5").Value = 10
$ws.Range("AZ5").Value = 151
$ws.Range("BA5").Value = 201
$ws.Range("G6").Value = 3.25
$ws.Range("I6").Value = 2.45
$ws.Range("K6").Value = 1.83
$ws.Range("O6").Value = 1.57
$ws.Range("P6").Value = 2.25
$ws.Range("Q6").Value = 2.88
$ws.Range("R6").Value = 1.4
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.57
$ws.Range("X6").Value = 15
$ws.Range("AB6").Value = 51
$ws.Range("AC6").Value = 5.5
$ws.Range("AN6").Value = 5
$ws.Range("AO6").Value = 21
$ws.Range("AR6").Value = 126
$ws.Range("G7").Value = 2.3
$ws.Range("H7").Value = 3.25
$ws.Range("I7").Value = 3
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.75
$ws.Range("S7").Value = 1.5
$ws.Range("T7").Value = 2.5
$ws.Range("U7").Value = 1.91
$ws.Range("V7").Value = 1.8
$ws.Range("Y7").Value = 9.5
$ws.Range("Z7").Value = 21
$ws.Range("AD7").Value = 6.5
$ws.Range("AE7").Value = 17
$ws.Range("AI7").Value = 12
$ws.Range("AT7").Value = 2.5
$ws.Range("G8").Value = 1.85
$ws.Range("I8").Value = 4.33
$ws.Range("J8").Value = 2.6
$ws.Range("L8").Value = 4.75
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9
$ws.Range("W8").Value = 6.5
$ws.Range("AD8").Value = 6
$ws.Range("AH8").Value = 21
$ws.Range("AI8").Value = 15
$ws.Range("AU8").Value = 8.5
$ws.Range("AX8").Value = 23
